$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 2173
$ws1.Range("F15").Value = 568
$ws1.Range("F16").Value = 414
$ws1.Range("F17").Value = 414
$ws1.Range("F20").Value = 2996
$ws1.Range("F23").Value = 3225
$ws1.Range("F28").Value = 739
$ws1.Range("F30").Value = 778
$ws1.Range("F31").Value = 759

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 69
$ws2.Range("F20").Value = 204
$ws2.Range("F21").Value = 145

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2940
$ws3.Range("F5").Value = 256
$ws3.Range("F6").Value = 411

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 256
$ws4.Range("F13").Value = 411
$ws4.Range("F14").Value = 2173
$ws4.Range("F18").Value = 69
$ws4.Range("F28").Value = 568
$ws4.Range("F29").Value = 414
$ws4.Range("F30").Value = 414
$ws4.Range("F35").Value = 2996
$ws4.Range("F37").Value = 3225
$ws4.Range("F45").Value = 204
$ws4.Range("F46").Value = 145
$ws4.Range("F48").Value = 739
$ws4.Range("F50").Value = 778
$ws4.Range("F51").Value = 759
